$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# The "query" column (B) for the ParticipantsTab row (row 2) was rewritten
# with an updated Cypher query that also pulls genomic_info and sorts the
# collected samples via apoc.coll.sort before joining them.
$newParticipantsQuery = "MATCH (p:participant)-->(s:study)`nOPTIONAL MATCH (samp:sample)-->(p)`nOPTIONAL MATCH (p)<--(diag:diagnosis)`nOPTIONAL MATCH (samp)<--(f:file)`nOPTIONAL MATCH (f)<--(g:genomic_info)`nWITH s, p, samp, f, g, diag`nWHERE f.file_type in ['FASTQ']`nwith p`nOPTIONAL MATCH (p)-->(s:study)`nOPTIONAL MATCH (samp:sample)-->(p)`nWITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp`nRETURN`ncoalesce(p.participant_id,'') as ``Participant ID``,`ncoalesce(s.study_name, '') as ``Study Name``,`ncoalesce(s.phs_accession,'') as ``Accession``,`ncoalesce(p.gender,'') as ``Gender``,`ncoalesce(apoc.text.join(samp, ','), '') as ``Samples```nORDER BY p.participant_id LIMIT 100"

$ws.Range("B2").Value = $newParticipantsQuery

# The cell wraps text, so the extra lines of the new query push row 2's
# height up from 186 to 279 points.
$ws.Rows.Item(2).RowHeight = 279

# Restore the view: scroll so row 3 is at the top-left and select B5 (matches
# the saved sheetView state in the edited file).
$aw = $excel.ActiveWindow
if ($aw) {
    $aw.ScrollRow = 3
    $aw.ScrollColumn = 1
}
$ws.Range("B5").Select() | Out-Null
